$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new data rows right after the header (new rows 2:3),
#    pushing the existing data/summary rows down by 2.
# ------------------------------------------------------------------
$ws.Rows("2:3").Insert()

# New row 2: "<0.25FLEXI"
$ws.Range("A2").Value = "<0.25FLEXI"
$ws.Range("B2").Value = 6821
$ws.Range("C2").Value = 6583
$ws.Range("D2").Value = 1667
$ws.Range("E2").Value = 2270
$ws.Range("F2").Value = 1168
$ws.Range("G2").Value = 516
$ws.Range("I2:M2").NumberFormat = "0.0%"
$ws.Range("I2").Formula = "=C2/(`$B2)"
$ws.Range("J2").Formula = "=D2/(`$B2)"
$ws.Range("K2").Formula = "=E2/(`$B2)"
$ws.Range("L2").Formula = "=F2/(`$B2)"
$ws.Range("M2").Formula = "=G2/(`$B2)"

# New row 3: "<0.50FLEXI"
$ws.Range("A3").Value = "<0.50FLEXI"
$ws.Range("B3").Value = 7731
$ws.Range("C3").Value = 7468
$ws.Range("D3").Value = 1917
$ws.Range("E3").Value = 2667
$ws.Range("F3").Value = 1377
$ws.Range("G3").Value = 589
$ws.Range("I3:M3").NumberFormat = "0.0%"
$ws.Range("I3").Formula = "=C3/(`$B3)"
$ws.Range("J3").Formula = "=D3/(`$B3)"
$ws.Range("K3").Formula = "=E3/(`$B3)"
$ws.Range("L3").Formula = "=F3/(`$B3)"
$ws.Range("M3").Formula = "=G3/(`$B3)"

# ------------------------------------------------------------------
# 2) Insert two new summary rows above the old "summary" rows (which,
#    after the shift above, start at row 17), pushing those down by 2.
# ------------------------------------------------------------------
$ws.Rows("17:18").Insert()

# New row 17: summary for "<0.25FLEXI" (references row 2)
$ws.Range("A17").Value = "<0.25FLEXI"
$ws.Range("B17:F17").HorizontalAlignment = -4108
$ws.Range("B17").Formula = '=C2&" ("&ROUND(100*I2,1)&"%)"'
$ws.Range("C17").Formula = '=E2&" ("&ROUND(100*K2,1)&"%)"'
$ws.Range("D17").Formula = '=F2&" ("&ROUND(100*L2,1)&"%)"'
$ws.Range("E17").Formula = '=D2&" ("&ROUND(100*J2,1)&"%)"'
$ws.Range("F17").Formula = '=G2&" ("&ROUND(100*M2,1)&"%)"'

# New row 18: summary for "<0.50FLEXI" (references row 3)
$ws.Range("A18").Value = "<0.50FLEXI"
$ws.Range("B18:F18").HorizontalAlignment = -4108
$ws.Range("B18").Formula = '=C3&" ("&ROUND(100*I3,1)&"%)"'
$ws.Range("C18").Formula = '=E3&" ("&ROUND(100*K3,1)&"%)"'
$ws.Range("D18").Formula = '=F3&" ("&ROUND(100*L3,1)&"%)"'
$ws.Range("E18").Formula = '=D3&" ("&ROUND(100*J3,1)&"%)"'
$ws.Range("F18").Formula = '=G3&" ("&ROUND(100*M3,1)&"%)"'

# ------------------------------------------------------------------
# 3) Selection / active cell moves to K21 (matches the saved sheetView).
# ------------------------------------------------------------------
$ws.Range("K21").Select()
